# Add 2022-Q3 data:
#  - insert a new "2022-Q3" worksheet right after the "总计" sheet, holding the
#    per-fund breakdown for the new quarter
#  - update the "总计" (summary) sheet: shift the existing two data rows down
#    by one and add a new top data row for 2022-Q3

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) "总计" sheet: push 2022-Q2 / 2021-Q4 rows down one row and add 2022-Q3
# ---------------------------------------------------------------------------

# Row 4 <- old row 3 (2021-Q4)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.23

# Row 3 <- old row 2 (2022-Q2)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.34

# Row 2 <- new 2022-Q3 summary
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.22

# Column A carries the bold/bordered/centered header-style (style index 2 in
# the original file); A2 already has it, propagate it onto the newly
# populated A4 (A3 already had it from before).
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, positioned right after "总计"
# ---------------------------------------------------------------------------

$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Columns that hold numeric-looking text (fund codes / percentages stored as
# text, matching the other quarter sheets) must be forced to Text format
# before the value is written, otherwise Excel auto-converts "008928" etc.
# into numbers and drops the leading zero / changes precision.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "008928"
$q3.Range("C2").Value = "泰达宏利中证主要消费红利指数A"
$q3.Range("D2").Value = "3.45"
$q3.Range("E2").Value = "93.07"
$q3.Range("F2").Value = "4.13"
$q3.Range("G2").Value = "0.1425"
$q3.Range("H2").Value = 7

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "008929"
$q3.Range("C3").Value = "泰达宏利中证主要消费红利指数C"
$q3.Range("D3").Value = "1.69"
$q3.Range("E3").Value = "93.07"
$q3.Range("F3").Value = "4.13"
$q3.Range("G3").Value = "0.0698"
$q3.Range("H3").Value = 7

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "501089"
$q3.Range("C4").Value = "方正富邦消费红利指数增强（LOF）"
$q3.Range("D4").Value = "0.31"
$q3.Range("E4").Value = "45.73"
$q3.Range("F4").Value = "2.01"
$q3.Range("G4").Value = "0.0062"
$q3.Range("H4").Value = 7

# The Text-format trick above minted a custom cell style; reset those cells
# back to the workbook's plain default style (no explicit `s`, like the
# other quarter sheets' data cells) by pasting formatting from a never-used,
# unstyled cell.
$q3.Range("Z100").Copy()
$q3.Range("B2:B4").PasteSpecial(-4122)
$q3.Range("D2:G4").PasteSpecial(-4122)

# Apply the same bold/bordered/centered style used on the other sheets'
# header rows and index (A) columns.
$totalSheet.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$q3.Range("A2:A4").PasteSpecial(-4122)
